# Updated parser to use TokenIteratorFieldRewriterSplit.
#
# The field delimiters "{" and "}" used to live in the same run as the
# first/last bit of token text ("{m" and "-&gt;sep(', ')}"). The new
# tokenizing field rewriter now always emits the delimiter in its own
# run, so we split those two runs here:
#   "{m"              -> "{"  + "m"
#   "-&gt;sep(', ')}" -> "-&gt;sep(', ')" + "}"

$d = $word.ActiveDocument

# The field is on the second paragraph of the document.
$fieldPara = $d.Paragraphs(2)

# --- Split "{m" into "{" and "m" -------------------------------------
$paraStart = $fieldPara.Range.Start
$openBrace = $d.Range($paraStart, $paraStart + 1)
# Toggling a character formatting property and reverting it is enough to
# force Word to materialize this single character as its own run,
# without altering the visible formatting of the text.
$openBrace.Font.Bold = 1
$openBrace.Font.Bold = 0

# --- Split "-&gt;sep(', ')}" into "-&gt;sep(', ')" and "}" -----------
# Re-fetch the paragraph end since the document length hasn't changed
# (no characters were inserted or removed above), but do it for safety.
$paraEnd = $fieldPara.Range.End
# paraEnd points just after the paragraph mark, so the closing brace is
# two positions before it.
$closeBracePos = $paraEnd - 2
$closeBrace = $d.Range($closeBracePos, $closeBracePos + 1)
$closeBrace.Font.Bold = 1
$closeBrace.Font.Bold = 0
